# Updates cryptos list values per the recorded price/volume diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.780.14"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "'3.407.12"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D5").Value = "'412.19"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'129.56"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.723"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("D11").Value = "'42.63"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "'0.0000218"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").Value = "'9.13"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").Value = "'3.949.89"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'20.39"
$ws.Range("D17").Value = "'3.403.36"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "'12.42"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").Value = "'1.07"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "'61.831.59"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'479.12"
$ws.Range("E21").Value = "  +16.89%  "
$ws.Range("D22").Value = "'90.71"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "'3.26"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").Value = "'13.12"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "'3.30"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("E26").Value = "  +10.41%  "
$ws.Range("D27").Value = "'33.13"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "'4.75"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'7.74"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.65"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'11.85"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "'0.167"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("D34").Value = "'40.97"
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("D36").Value = "'57.09"
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("D40").Value = "'149.04"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  +4.53%  "
$ws.Range("D45").Value = "'2.59"
$ws.Range("E45").Value = "  +7.38%  "
$ws.Range("D46").Value = "'4.19"
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("E47").Value = "  +19.02%  "
$ws.Range("D48").Value = "'16.44"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "'0.0₃0531"
$ws.Range("E49").Value = "  +17.13%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'22.09"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "'112.88"
$ws.Range("E51").Value = "  +14.27%  "
